# Adds an NSA-funding acknowledgement paragraph after the trademark
# paragraph near the end of the document, moves the "_GoBack" bookmark
# down so it still sits in its own (now later) empty paragraph, and
# refreshes the cached PAGE field result in the default footer so the
# printed page count stays in sync with the newly added content.

$d = $word.ActiveDocument

$t1 = "Project sponsored by the National Security Agency under grant Number H98230-17-1-0199."
$t2 = " "
$t3 = "The United States Government is authorized to reproduce and distribute reprints notwithstanding any copyright notation herein."

# Step 1: Insert a new centered paragraph right after the trademark
# paragraph, holding the full acknowledgement text (still as one run).
$null = $d.Content.Find.Execute("respective holders.", $true, $false, $false, $false, $false, $true, 1, $false, "respective holders.^p" + $t1 + $t2 + $t3, 2)

# Step 2: Split off a second, empty, centered paragraph right after it
# (this is where the "_GoBack" bookmark will live).
$null = $d.Content.Find.Execute($t3, $true, $false, $false, $false, $false, $true, 1, $false, $t3 + "^p", 2)

# Step 3: Apply the 24-half-point (12pt) size to the first and third
# runs of the acknowledgement paragraph, leaving the middle space
# run with no explicit size - this naturally splits the paragraph
# into the three runs seen in the final document.
$p = $d.Paragraphs.Item(99)
$pStart = $p.Range.Start

$r1Start = $pStart
$r1End = $r1Start + $t1.Length
$r2Start = $r1End
$r2End = $r2Start + $t2.Length
$r3Start = $r2End
$r3End = $r3Start + $t3.Length

$run1 = $d.Range($r1Start, $r1End)
$run1.Font.Size = 12

$run3 = $d.Range($r3Start, $r3End)
$run3.Font.Size = 12

# Step 4: Move the "_GoBack" bookmark into the new trailing empty
# paragraph (Bookmarks.Add replaces any existing bookmark of the same
# name, so the old location is cleared automatically).
$goBackPara = $d.Paragraphs.Item(100)
$d.Bookmarks.Add("_GoBack", $goBackPara.Range)

# Step 5: Refresh the cached page-number field result in the default
# footer so it reflects the document's new pagination.
$footer = $d.Sections.Item(1).Footers.Item(1)
$pageField = $footer.Range.Fields.Item(1)
$null = $pageField.Result.Find.Execute("2", $true, $false, $false, $false, $false, $true, 1, $false, "5", 2)
